$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.882.89"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "1.669.21"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").Value = "'215.62"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("E6").Value = "  +4.32%  "

$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  +1.85%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("D10").Value = "'20.28"
$ws.Range("E10").Value = "  +3.73%  "

$ws.Range("D11").Value = "'0.0893"
$ws.Range("E11").Value = "  +3.74%  "

$ws.Range("D12").Value = "1.905.20"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").Value = "1.672.61"
$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("E15").Value = "  +1.16%  "

$ws.Range("D16").Value = "'65.74"
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").Value = "26.902.30"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "'231.92"
$ws.Range("E18").Value = "  -4.56%  "

$ws.Range("D19").Value = "'7.78"
$ws.Range("E19").Value = "  -1.24%  "

$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'2.21"
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").Value = "'9.19"
$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("D25").Value = "'145.53"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").Value = "'7.14"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("D28").Value = "'15.92"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").Value = "1.464.36"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("D34").Value = "'3.16"
$ws.Range("E34").Value = "  +3.63%  "

$ws.Range("E35").Value = "  +3.60%  "

$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "'0.900"
$ws.Range("E37").Value = "  +0.68%  "

$ws.Range("D38").Value = "'0.571"
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").Value = "'5.84"
$ws.Range("E40").Value = "  -2.37%  "

$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").Value = "'2.29"
$ws.Range("E42").Value = "  +2.50%  "

$ws.Range("D43").Value = "'65.76"
$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("E44").Value = "  +6.80%  "

$ws.Range("D45").Value = "1.814.54"
$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").Value = "'90.43"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0508"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.57"
$ws.Range("E51").Value = "  +0.37%  "

Write-Host "Applied cryptos update"